# Refresh the scraped crypto price/volume table (GitHub Actions sync).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # The sheet stores every Price/Volume cell as literal text (values like
    # "0.530" or "1.00" must keep their trailing zeros). Plain
    # `$range.Value = $value` lets Excel auto-convert simple decimals to a
    # Number, which would lose them -- so force a Text format first, then
    # clear the temporary format again so the cell keeps its original
    # (default) style and only the content changes.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

$ws.Range('D2').Value = '60.112.13'
$ws.Range('E2').Value = '  +2.49%  '
$ws.Range('D3').Value = '3.194.60'
$ws.Range('E3').Value = '  +1.39%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('E5').Value = '  +1.10%  '
Set-TextValue $ws.Range('D6') '145.43'
$ws.Range('E6').Value = '  +4.16%  '
$ws.Range('E7').Value = '  -0.01%  '
Set-TextValue $ws.Range('D8') '0.530'
$ws.Range('E8').Value = '  -0.68%  '
Set-TextValue $ws.Range('D9') '7.33'
$ws.Range('E9').Value = '  +0.20%  '
$ws.Range('E10').Value = '  +2.28%  '
Set-TextValue $ws.Range('D11') '0.431'
$ws.Range('E11').Value = '  -0.25%  '
$ws.Range('D12').Value = '3.748.73'
$ws.Range('E12').Value = '  +1.46%  '
Set-TextValue $ws.Range('D13') '0.137'
$ws.Range('E13').Value = '  -2.12%  '
$ws.Range('E14').Value = '  -0.39%  '
$ws.Range('E15').Value = '  +0.94%  '
$ws.Range('D16').Value = '60.151.92'
$ws.Range('E16').Value = '  +2.49%  '
$ws.Range('D17').Value = '3.220.78'
$ws.Range('E17').Value = '  +2.40%  '
$ws.Range('E18').Value = '  +0.32%  '
$ws.Range('E19').Value = '  +2.37%  '
Set-TextValue $ws.Range('D20') '8.23'
$ws.Range('E20').Value = '  +0.99%  '
Set-TextValue $ws.Range('D21') '371.60'
$ws.Range('E21').Value = '  -0.16%  '
Set-TextValue $ws.Range('D22') '1.00'
$ws.Range('E22').Value = '  +0.03%  '
Set-TextValue $ws.Range('D23') '0.522'
$ws.Range('E23').Value = '  -0.10%  '
Set-TextValue $ws.Range('D24') '69.51'
$ws.Range('E24').Value = '  -0.40%  '
$ws.Range('E25').Value = '  +1.27%  '
Set-TextValue $ws.Range('D26') '8.62'
$ws.Range('E26').Value = '  +4.39%  '
$ws.Range('E27').Value = '  -0.09%  '
$ws.Range('E28').Value = '  +1.91%  '
Set-TextValue $ws.Range('D29') '22.51'
$ws.Range('E29').Value = '  +1.80%  '
$ws.Range('E30').Value = '  +0.91%  '
$ws.Range('E31').Value = '  +0.75%  '
Set-TextValue $ws.Range('D32') '5.30'
$ws.Range('E32').Value = '  +2.51%  '
$ws.Range('E33').Value = '  +2.64%  '
Set-TextValue $ws.Range('D34') '6.56'
$ws.Range('E34').Value = '  +4.66%  '
Set-TextValue $ws.Range('D35') '156.80'
$ws.Range('E35').Value = '  -1.39%  '
$ws.Range('E36').Value = '  +1.88%  '
$ws.Range('B37').Value = 'EnergySwap'
$ws.Range('C37').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range('D37') '26.37'
$ws.Range('E37').Value = '  +4.76%  '
$ws.Range('B38').Value = 'Maker'
$ws.Range('C38').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D38').Value = '2.823.74'
$ws.Range('E38').Value = '  +7.19%  '
Set-TextValue $ws.Range('D39') '0.0705'
$ws.Range('E39').Value = '  +2.83%  '
Set-TextValue $ws.Range('D40') '0.0309'
$ws.Range('E40').Value = '  +8.16%  '
Set-TextValue $ws.Range('D41') '1.68'
$ws.Range('E41').Value = '  +0.19%  '
$ws.Range('E42').Value = '  -0.11%  '
Set-TextValue $ws.Range('D43') '39.97'
$ws.Range('E43').Value = '  +2.72%  '
$ws.Range('E44').Value = '  +1.62%  '
$ws.Range('E45').Value = '  +0.67%  '
$ws.Range('D46').Value = '3.237.87'
$ws.Range('E46').Value = '  +1.34%  '
$ws.Range('E47').Value = '  +0.70%  '
Set-TextValue $ws.Range('D48') '6.15'
$ws.Range('E48').Value = '  -0.73%  '
Set-TextValue $ws.Range('D49') '20.73'
$ws.Range('E49').Value = '  +1.95%  '
Set-TextValue $ws.Range('D50') '0.795'
$ws.Range('E50').Value = '  +5.22%  '
$ws.Range('E51').Value = '  +0.07%  '
